$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row
$ws.Range("A1").Value = "Toimittajanro"
$ws.Range("B1").Value = "Y-Tunnus"
$ws.Range("C1").Value = "Toimittaja"
$ws.Range("D1").Value = "Summa"
$ws.Range("E1").Value = "Osoite"
$ws.Range("F1").Value = "Postinumero"
$ws.Range("G1").Value = "Kaupunki"
$ws.Range("H1").Value = "Tilinro"
$ws.Range("I1").Value = "Pankki"

# Data rows
$ws.Range("A2").Value = 11
$ws.Range("B2").Value = "0725937-3"
$ws.Range("C2").Value = "E-K SOSIAALI- JA TERVEYSPIIRI"
$ws.Range("D2").Value = 21509717.18

$ws.Range("A3").Value = 73978
$ws.Range("B3").Value = "2977551-2"
$ws.Range("C3").Value = "PEAB INDUSTRI OY"
$ws.Range("D3").Value = 860908.83

$ws.Range("A4").Value = 77161
$ws.Range("B4").Value = "2752472-8"
$ws.Range("C4").Value = "PALLAS RAKENNUS KAAKKOIS-SUOMI"
$ws.Range("D4").Value = 782500

$ws.Range("A5").Value = 38975
$ws.Range("B5").Value = "0162017-2"
$ws.Range("C5").Value = "SAIMAAN TUKIPALVELUT OY"
$ws.Range("D5").Value = 560267.46

$ws.Range("A6").Value = 67122
$ws.Range("B6").Value = "1565583-5"
$ws.Range("C6").Value = "YIT SUOMI OY"
$ws.Range("D6").Value = 400000
